# Recolor the "intervention_type" legend from the black/red/orange/green
# palette to a blue/red/orange/green palette:
#   ⬛ -> 📘   🟥 -> 📕   🟧 -> 📙   🟩 -> 📗   noir -> bleu
#
# Column A holds the colored-square glyph ("statut"), column B holds the
# matching French color label ("statut_label"). Walk every used row and
# rewrite any cell that still carries one of the old values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valueA = $cellA.Value()

    if ($valueA -eq "⬛") { $cellA.Value = "📘" }
    elseif ($valueA -eq "🟥") { $cellA.Value = "📕" }
    elseif ($valueA -eq "🟧") { $cellA.Value = "📙" }
    elseif ($valueA -eq "🟩") { $cellA.Value = "📗" }

    $cellB = $ws.Cells.Item($r, 2)
    $valueB = $cellB.Value()

    if ($valueB -eq "noir") { $cellB.Value = "bleu" }
}
